$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-PlainValue($cellAddr, $val) {
    $ws.Range($cellAddr).Value = $val
}

# --- Rows 44 & 45 swap: VeChain/Maker order changes with updated data ---
Set-PlainValue "B44" "Maker"
Set-PlainValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-PlainValue "D44" "2.692.52"
Set-PlainValue "E44" "  +0.85%  "

Set-PlainValue "B45" "VeChain"
Set-PlainValue "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0342"
Set-PlainValue "E45" "  -1.58%  "

# --- Other row 2-51 simple cell value updates (price & volume %) ---
Set-PlainValue "D2" "60.911.07"
Set-PlainValue "E2" "  +0.22%  "
Set-PlainValue "D3" "2.918.60"
Set-PlainValue "E3" "  +0.22%  "
Set-PlainValue "E4" "  +0.02%  "
Set-TextValue "D5" "590.11"
Set-PlainValue "E5" "  +1.10%  "
Set-TextValue "D6" "146.41"
Set-PlainValue "E6" "  +1.73%  "
Set-PlainValue "E7" "  +0.01%  "
Set-PlainValue "E8" "  +0.80%  "
Set-PlainValue "E9" "  +1.12%  "
Set-PlainValue "E10" "  -0.34%  "
Set-TextValue "D11" "0.440"
Set-PlainValue "E11" "  -1.36%  "
Set-PlainValue "E12" "  +0.08%  "
Set-TextValue "D13" "33.54"
Set-PlainValue "E13" "  +0.00%  "
Set-PlainValue "D15" "3.401.18"
Set-PlainValue "E15" "  +0.17%  "
Set-PlainValue "D16" "60.818.63"
Set-PlainValue "E16" "  +0.14%  "
Set-PlainValue "D18" "2.916.79"
Set-PlainValue "E18" "  +0.21%  "
Set-TextValue "D19" "430.46"
Set-PlainValue "E19" "  +0.11%  "
Set-PlainValue "E20" "  -1.80%  "
Set-PlainValue "E21" "  -0.59%  "
Set-TextValue "D22" "7.06"
Set-PlainValue "E22" "  -0.86%  "
Set-TextValue "D23" "81.36"
Set-PlainValue "E23" "  +1.26%  "
Set-TextValue "D24" "10.94"
Set-PlainValue "E24" "  +1.45%  "
Set-PlainValue "E25" "  -0.42%  "
Set-TextValue "D26" "11.84"
Set-PlainValue "E26" "  -0.35%  "
Set-PlainValue "E27" "  +0.05%  "
Set-TextValue "D28" "2.26"
Set-PlainValue "E28" "  +4.54%  "
Set-PlainValue "E29" "  +0.11%  "
Set-PlainValue "E30" "  -2.72%  "
Set-TextValue "D31" "26.62"
Set-PlainValue "E31" "  +0.70%  "
Set-PlainValue "E32" "  +1.84%  "
Set-PlainValue "E33" "  +0.02%  "
Set-PlainValue "D34" "0.0₃0857"
Set-PlainValue "E34" "  -0.87%  "
Set-PlainValue "E35" "  +0.13%  "
Set-PlainValue "E36" "  -0.56%  "
Set-PlainValue "E37" "  +0.61%  "
Set-PlainValue "E38" "  -1.13%  "
Set-PlainValue "E39" "  -3.50%  "
Set-PlainValue "E40" "  -1.17%  "
Set-PlainValue "E41" "  -4.14%  "
Set-TextValue "D42" "40.20"
Set-PlainValue "E42" "  -2.71%  "
Set-TextValue "D43" "379.72"
Set-PlainValue "E43" "  +1.74%  "
Set-TextValue "D46" "133.40"
Set-PlainValue "E46" "  +1.14%  "
Set-PlainValue "E47" "  -0.03%  "
Set-PlainValue "E48" "  -2.13%  "
Set-PlainValue "E49" "  -0.59%  "
Set-PlainValue "E50" "  -3.64%  "
Set-PlainValue "E51" "  -0.13%  "

Write-Host "Edit complete"
